$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the requisito texts between row 23 and row 24 so that the
# "LOT2028 - Tecnologia de Processos Fermentativos (Requisito fraco)" entry
# comes before the "LOT2038 - Tecnologia de Bebidas (Indicação de Conjunto)" entry.
$lot2038 = "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)`n"
$lot2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"

$ws.Range("B23").Value = $lot2028
$ws.Range("C23").Value = $lot2028

$ws.Range("B24").Value = $lot2038
$ws.Range("C24").Value = $lot2038
